$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 17.8138221557582
$ws.Range("F2").Value = 0.5695516942056902
$ws.Range("K2").Value = 0.7364114393528272

$ws.Range("E3").Value = 6.941370064669262
$ws.Range("F3").Value = 0.3563056214683741
$ws.Range("G3").Value = -13.21
$ws.Range("K3").Value = 0.8706318161174728

$ws.Range("E4").Value = 20.06614449729091
$ws.Range("F4").Value = 0.7595227110328777
$ws.Range("G4").Value = -13.21
$ws.Range("K4").Value = 0.4700464414683252

$ws.Range("E5").Value = 22.34089329233793
$ws.Range("F5").Value = 0.845624122889358
$ws.Range("G5").Value = -13.21
$ws.Range("K5").Value = 0.6539970279579241

$ws.Range("E6").Value = 20.06614449729091
$ws.Range("F6").Value = 0.7595227110328777
$ws.Range("G6").Value = -13.21
$ws.Range("K6").Value = 0.4700464414683252

$ws.Range("E7").Value = 22.34089329233793
$ws.Range("F7").Value = 0.845624122889358
$ws.Range("G7").Value = -13.21
$ws.Range("K7").Value = 0.6539970279579241

$ws.Range("E8").Value = 9.055951952863756
$ws.Range("F8").Value = 0.5217370667757679
$ws.Range("G8").Value = -13.21
$ws.Range("K8").Value = 0.5944602021812885

$ws.Range("E9").Value = 2.844697658355452
$ws.Range("F9").Value = 0.2800995812387078
$ws.Range("K9").Value = 0.4965408554272375

$ws.Range("E10").Value = 12.89443785194213
$ws.Range("F10").Value = 0.4882834990804356
$ws.Range("G10").Value = -13.21
$ws.Range("K10").Value = 0.9757307980375048

$ws.Range("E11").Value = 11.19787297467833
$ws.Range("F11").Value = 0.3487541619380949
$ws.Range("G11").Value = -13.21
$ws.Range("K11").Value = 0.8553152162217491

$ws.Range("E12").Value = 12.48937108880223
$ws.Range("F12").Value = 0.7667243171293792
$ws.Range("K12").Value = 0.6666701490697698

$ws.Range("E13").Value = 7.096204240678009
$ws.Range("F13").Value = 0.4864762432011431
$ws.Range("G13").Value = -13.21
$ws.Range("K13").Value = 0.7490918652199673

$ws.Range("E14").Value = 2.761960143974892
$ws.Range("F14").Value = 0.2975663924620373
$ws.Range("G14").Value = -13.21
$ws.Range("K14").Value = 0.4787472366890582

$ws.Range("E15").Value = 2.677177424089021
$ws.Range("F15").Value = 0.2727280963681851
$ws.Range("G15").Value = -13.21
$ws.Range("K15").Value = 0.5023728944354221

$ws.Range("E16").Value = 12.52385416926674
$ws.Range("F16").Value = 0.3900514209184453
$ws.Range("G16").Value = -13.21

$ws.Range("E17").Value = 24.24659203144867
$ws.Range("F17").Value = 0.9074665159154658
$ws.Range("G17").Value = -13.21
$ws.Range("K17").Value = 0.9970639298715548

$ws.Range("E18").Value = 10.42241535579358
$ws.Range("F18").Value = 0.6135132742779132
$ws.Range("G18").Value = -13.21
$ws.Range("K18").Value = 0.6241356408749192

$ws.Range("E19").Value = 3.969093041339537
$ws.Range("F19").Value = 0.4540541189076495
$ws.Range("K19").Value = 0.4861187771148835

$ws.Range("E20").Value = 0.09996363129261256
$ws.Range("F20").Value = 0.2355252494837661
$ws.Range("K20").Value = 0.9798166665928864

$ws.Range("E21").Value = 0.1141673728959025
$ws.Range("F21").Value = 0.2428831806141352
$ws.Range("K21").Value = 0.9759764580263021

$ws.Range("E22").Value = 12.48937108880226
$ws.Range("F22").Value = 0.5410628137990361
$ws.Range("G22").Value = -13.21
$ws.Range("K22").Value = 0.6612464686765634

$ws.Range("E23").Value = 2.511143878215948
$ws.Range("F23").Value = 0.3475710879338727
$ws.Range("K23").Value = 0.4645429078389048

$ws.Range("E24").Value = 3.969093041339538
$ws.Range("F24").Value = 0.3806549572065452
$ws.Range("G24").Value = -13.21
$ws.Range("K24").Value = 0.4861187771148942

$ws.Range("E25").Value = 2.957467104488082
$ws.Range("F25").Value = 0.306078473194127
$ws.Range("G25").Value = -13.21
$ws.Range("K25").Value = 0.4724524112256044

$ws.Range("E26").Value = 9.725652207600666
$ws.Range("F26").Value = 0.7116896754722978
$ws.Range("K26").Value = 0.669206989419655

$ws.Range("E27").Value = 0.1013986989748657
$ws.Range("F27").Value = 0.2406859896067657
$ws.Range("K27").Value = 0.9800385302110093

$ws.Range("E28").Value = 6.766253892867912
$ws.Range("F28").Value = 0.2961120057049863
$ws.Range("G28").Value = -13.21
$ws.Range("K28").Value = 0.7980874276204809

$ws.Range("E29").Value = 0.1308423263314998
$ws.Range("F29").Value = 0.2719810757756068
$ws.Range("K29").Value = 0.9575254300092648

$ws.Range("E30").Value = 0.007105612876140991
$ws.Range("F30").Value = 0.06718414792941983
